# edit.ps1 - applies the "create default templates and update about screen" commit
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Abstract paragraph (paragraph 4) text tweaks
# ---------------------------------------------------------------------------

$d.Content.Find.Execute(
    'titled ‘Standardizing Nomenclatures in Radiation Oncology’. Unfortunately, despite having increased guidance on nomenclature, the burden of converting ',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'titled ‘Standardizing Nomenclatures in Radiation Oncology’ to assist in this nomenclature. Unfortunately, the burden of converting ',
    2) | Out-Null

$d.Content.Find.Execute(
    'the treatment planning system used. Our work',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'the treatment planning system implemented. Our work',
    2) | Out-Null

$d.Content.Find.Execute(
    'windows system and',
    $true, $false, $false, $false, $false, $true, 1, $false,
    'Windows system, and has',
    2) | Out-Null

Write-Output 'abstract text updated'

# ---------------------------------------------------------------------------
# 2) Replace the trailing empty paragraph with the new "Introduction" section
#    and the start of "Methods"
# ---------------------------------------------------------------------------

$pIntroHeading = $d.Paragraphs.Item(5)
$pIntroHeading.Range.Text = 'Introduction'
$pIntroHeading.Style = 'Heading1'

$pIntro1 = $d.Paragraphs.Add($pIntroHeading.Range)
$pIntro1.Style = 'Normal'
$pIntro1.Range.Text = 'In the creation of a treatment plan within radiation oncology, regions of interest (ROIs) must be defined. These ROIs can be the target of radiation therapy, organs at risk (OARs), or contrast agents, etc. While the Digital Imaging and Communications in Medicine (DICOM) provides a standard for communicating these generated structures (RT-Structures) in treatment planning systems, the creation of the RT-Structures is often left to the treatment planning system.'

$pIntro2 = $d.Paragraphs.Add($pIntro1.Range)
$pIntro2.Style = 'Normal'
$pIntro2.Range.Text = 'The most important ROIs will vary based on the treatment site. For example, when treating disease in the skull, the Liver will likely not be of interest. Depending on the treatment planning system, the user will then be required to manually create each ROI, individually labeling the structures involved (‘Brain’, ‘Brainstem’, etc.). This can be not only tedious, but also error prone (‘Brian’ instead of ‘Brain’). Furthermore, the naming of an ROI can vary from person to person (‘Lung_R’ vs ‘Right Lung’). While several treatment planning systems provide a method of creating templates to automatically create the desired ROIs based on the treatment site[ref for varian, others?], these templates must be created manually.'

$pIntro3 = $d.Paragraphs.Add($pIntro2.Range)
$pIntro3.Style = 'Normal'
$pIntro3.Range.Text = 'The American Association of Physics in Medicine (AAPM) has created Report 263 titled ‘Standardizing Nomenclatures in Radiation Oncology’, whose purpose is to provide guidance on naming of ROIs. Unfortunately, adoption of this can be difficult based on the tools available in the clinic. In a recent survey provided by TG-263, 689 responses from members of AAPM, the American Society for Radiation Oncology (ASTRO), and the American Association of Medical Dosimetrists (AAMD) were asked about their likelihood for adopting TG-263. For respondents who had not yet adopted the new nomenclature, the majority stated that the largest hurdle was difficulty with retraining staff and/or a lack of time/resources to create new templates. With this work, we hope to provide a simple, server based system that will automatically create the desired RT-Structure files, and provide several ‘standard’ templates for commonly treated sites.'

$pMethodsHeading = $d.Paragraphs.Add($pIntro3.Range)
$pMethodsHeading.Style = 'Heading1'
$pMethodsHeading.Range.Text = 'Methods'

$pMethods1 = $d.Paragraphs.Add($pMethodsHeading.Range)
$pMethods1.Style = 'Normal'
$pMethods1.Range.Text = 'The program '

Write-Output 'introduction/methods sections added'
Write-Output "paragraph count: $($d.Paragraphs.Count)"
